# Auto update: 2025-12-06 21:20:02
# Corrects the ticker/name pairing for rows 3-5 and refreshes the computed
# score columns (K, N) for all data rows on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2 (HDKSOE / 009540.KS) : only final score columns change ---
$ws.Range("K2").Value = 59.7
$ws.Range("N2").Value = 52.28493729186943

# --- Row 3 : now Hanwha Ocean / 042660.KS ---
$ws.Range("B3").Value = "Hanwha Ocean"
$ws.Range("C3").Value = "042660.KS"
$ws.Range("D3").Value = 108500
$ws.Range("E3").Value = 21
$ws.Range("F3").Value = 0.65
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 70
$ws.Range("I3").Value = 80
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 53.7
$ws.Range("N3").Value = 52.28493729186943

# --- Row 4 : now SamsungHvyInd / 010140.KS ---
$ws.Range("B4").Value = "SamsungHvyInd"
$ws.Range("C4").Value = "010140.KS"
$ws.Range("D4").Value = 25450
$ws.Range("E4").Value = 44.4
$ws.Range("F4").Value = 3.46
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 60
$ws.Range("I4").Value = 70
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 52.7
$ws.Range("N4").Value = 52.28493729186943

# --- Row 5 : now HD HYUNDAI MIPO / 010620.KS ---
$ws.Range("B5").Value = "HD HYUNDAI MIPO"
$ws.Range("C5").Value = "010620.KS"
$ws.Range("D5").Value = 223000
$ws.Range("E5").Value = 26.8
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = 70
$ws.Range("J5").Value = 63
$ws.Range("K5").Value = 46.7
$ws.Range("N5").Value = 52.28493729186943
